# PNAD 2009 - furto - "correção nos dados e inicio da analise PNAD 2009"
#
# The sheet had a two-row header (row 1 "total" / row 2 a pandas
# multi-index artifact "unnamed: 1_level_1") and its data block was
# interleaved with sub-category header rows ("sexo", "cor ou raça",
# "grupos de idade", "nível de instrução", "classes de rendimento mensal
# domiciliar per capita") that carried no data of their own, plus two
# trailing footnote rows. This edit:
#   1. Fixes the stray "unnamed: 1_level_1" header cell to read "total".
#   2. Removes the label-only sub-category header rows and the two
#      footnote rows, letting the data rows shift up into a clean,
#      contiguous table (A1:G28).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the B2 header cell ---
$ws.Range("B2").Value = "total"

# --- 2. Delete the empty header/footnote rows, bottom-up so the row
#        numbers of rows still waiting to be deleted don't shift ---
$ws.Rows(35).EntireRow.Delete()   # "(1) inclusive as pessoas de cor ou raça amarela..."
$ws.Rows(34).EntireRow.Delete()   # "fonte: ibge, diretoria de pesquisas..."
$ws.Rows(27).EntireRow.Delete()   # "classes de rendimento mensal domiciliar per capita"
$ws.Rows(19).EntireRow.Delete()   # "nível de instrução"
$ws.Rows(13).EntireRow.Delete()   # "grupos de idade"
$ws.Rows(8).EntireRow.Delete()    # "cor ou raça"
$ws.Rows(5).EntireRow.Delete()    # "sexo"

Write-Output "done"
